# Apply the commit "A bit more complex still":
#  - Summary!B5 label changes from the old "Alpha, Jan" example to "Table search"
#    (the config block below it is being generalised into a table/regex search demo).
#  - Config sheet's second search block (rows 13-14) is changed from a literal
#    "Jan" search against "Alpha" to a regex "Mar.*" search (matches) against "Beta".
#  - The previously-active "Config" sheet's selection moves to C15, and "Summary"
#    becomes the active sheet with its selection on D8.
#  - Column widths are nudged (Excel's best-fit) to account for the new, longer text.

$wb = $excel.ActiveWorkbook
$sum = $wb.Worksheets.Item("Summary")
$cfg = $wb.Worksheets.Item("Config")

# --- Config sheet edits first (so new shared strings land in the same order
#     the authoring app produced: matches, Mar.*, Target intersection..., Table search) ---
$cfg.Range("C13").Value = "matches"
$cfg.Range("D13").Value = "Mar.*"
$cfg.Range("F13").Value = "Target intersection of Mar (using a regex search in caes it says March) and Beta"
$cfg.Range("D14").Value = "Beta"

# --- Summary sheet edit ---
$sum.Range("B5").Value = "Table search"

# --- Column width nudges (best-fit growth from the new, longer cell text) ---
$sum.Columns.Item(2).ColumnWidth = 10.33
$sum.Columns.Item(3).ColumnWidth = 11.67
$cfg.Columns.Item(2).ColumnWidth = 17.06
$cfg.Columns.Item(3).ColumnWidth = 6.83

# --- Selection / active sheet changes ---
$cfg.Activate()
$cfg.Range("C15").Select()

$sum.Activate()
$sum.Range("D8").Select()
